$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.553.36"
$ws.Range("D3").Value = "1.679.12"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'219.95"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'29.98"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").Value = "'0.0628"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "1.920.61"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").Value = "'10.42"
$ws.Range("E13").Value = "  +12.63%  "
$ws.Range("E14").Value = "  +9.73%  "
$ws.Range("D15").Value = "1.675.34"
$ws.Range("E15").Value = "  +2.94%  "
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "30.577.31"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "'66.38"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").Value = "'244.82"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'10.16"
$ws.Range("E22").Value = "  +3.60%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.28"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").Value = "'157.49"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'15.88"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("D28").Value = "'6.71"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'0.0498"
$ws.Range("E30").Value = "  +2.40%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("D32").Value = "'3.48"
$ws.Range("D33").Value = "1.509.28"
$ws.Range("E33").Value = "  +5.56%  "
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("D35").Value = "'1.75"
$ws.Range("E35").Value = "  +6.92%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "'83.61"
$ws.Range("E37").Value = "  +10.69%  "
$ws.Range("D38").Value = "'0.0179"
$ws.Range("E38").Value = "  +5.38%  "
$ws.Range("D39").Value = "'0.591"
$ws.Range("E39").Value = "  +7.39%  "
$ws.Range("D40").Value = "'2.72"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "'5.57"
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("D48").Value = "'51.55"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").Value = "1.813.55"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "'94.61"
$ws.Range("E50").Value = "  +6.28%  "
$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +2.71%  "
